# Generate Report for Handoff
# Adds a new file entry (d0e5da63-70f7-42fb-9eee-88ffc50b14cb) as row 7 to each
# of the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$fileGuid   = "d0e5da63-70f7-42fb-9eee-88ffc50b14cb"
$mdName     = $fileGuid + ".md"
$zhXlfName  = $fileGuid + ".b785921fcbf62b7f1a9196392702a6f75eaa19f1.zh-cn.xlf"
$deXlfName  = $fileGuid + ".b785921fcbf62b7f1a9196392702a6f75eaa19f1.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet (row 7)
# ---------------------------------------------------------------------------
$wsOverview.Range("A7").Value2 = $mdName
$wsOverview.Range("B7").Value2 = "Ready for handoff"
$wsOverview.Range("C7").Value2 = "Ready for handoff"
$wsOverview.Range("D7").Value2 = "2016-03-23 06:33:35"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A7"),
    "https://github.com/OpenLocalizationTest/oltest/blob/a9d5f8b4195760fee57e469e4059778a46370ecb/e2e/" + $mdName,
    [Type]::Missing,
    [Type]::Missing,
    $mdName
)

# ---------------------------------------------------------------------------
# zh-cn sheet (row 7)
# ---------------------------------------------------------------------------
$wsZhCn.Range("A7").Value2 = $mdName
$wsZhCn.Range("B7").Value2 = ".md"
$wsZhCn.Range("C7").Value2 = "Ready for handoff"
$wsZhCn.Range("D7").Value2 = $zhXlfName
$wsZhCn.Range("E7").Value2 = "2016-03-23 06:33:27"
$wsZhCn.Range("H7").Value2 = "0001-01-01 00:00:00"
$wsZhCn.Range("J7").Value2 = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A7"),
    "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/09b9d1bad2ca176fba1eced6f9ec3bc701eae891/e2e/" + $mdName,
    [Type]::Missing,
    [Type]::Missing,
    $mdName
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D7"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/08116dd94051fe66c672a1ec9c81645d6eaf5e5e/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/" + $zhXlfName,
    [Type]::Missing,
    [Type]::Missing,
    $zhXlfName
)

# ---------------------------------------------------------------------------
# de-de sheet (row 7)
# ---------------------------------------------------------------------------
$wsDeDe.Range("A7").Value2 = $mdName
$wsDeDe.Range("B7").Value2 = ".md"
$wsDeDe.Range("C7").Value2 = "Ready for handoff"
$wsDeDe.Range("D7").Value2 = $deXlfName
$wsDeDe.Range("E7").Value2 = "2016-03-23 06:33:35"
$wsDeDe.Range("H7").Value2 = "0001-01-01 00:00:00"
$wsDeDe.Range("J7").Value2 = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A7"),
    "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/e7a41e75ab426537c1cdc93e95a59b5ab64aea8b/e2e/" + $mdName,
    [Type]::Missing,
    [Type]::Missing,
    $mdName
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D7"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1640cc960c049c99c1017b7553ce8fa4e5006716/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/" + $deXlfName,
    [Type]::Missing,
    [Type]::Missing,
    $deXlfName
)
